$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numericTextCells = @("B6","D6","B7","D7","B8","D8")
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("A6").Value = "Ольга"
$ws.Range("B6").Value = "89865"
$ws.Range("C6").Value = "22.02.2023"
$ws.Range("D6").Value = "10"
$ws.Range("E6").Value = "Нет, сегодня впервые рисовала👍"
$ws.Range("F6").Value = "запомнили основные правила и технику безопасности для экологичного рисования"
$ws.Range("G6").Value = "Вопрос"
$ws.Range("H6").Value = "Воодушевление, уверенность, спокойствие"
$ws.Range("I6").Value = "Отзыв"
$ws.Range("J6").Value = "Да"
$ws.Range("K6").Value = "Нет"
$ws.Range("L6").Value = "Увеличение дохода"

$ws.Range("A7").Value = "Булат"
$ws.Range("B7").Value = "8965"
$ws.Range("C7").Value = "Мораль 22.03.2023"
$ws.Range("D7").Value = "8"
$ws.Range("E7").Value = "Рисовала на бесплатном марафоне😁"
$ws.Range("F7").Value = "запомнили основные правила и технику безопасности для экологичного рисования"
$ws.Range("G7").Value = "Всё понятно!"
$ws.Range("H7").Value = "Напряжение усилилось"
$ws.Range("I7").Value = "Нормас"
$ws.Range("J7").Value = "Да"
$ws.Range("K7").Value = "Да"
$ws.Range("L7").Value = "Здоровье"

$ws.Range("A8").Value = "Bulat"
$ws.Range("B8").Value = "11233"
$ws.Range("C8").Value = "Fhjj"
$ws.Range("D8").Value = "10"
$ws.Range("E8").Value = "Рисовала другой алгоритм на очном МК"
$ws.Range("F8").Value = "вошли в состояние медитации, появилось ощущение гармонии"
$ws.Range("G8").Value = "Всё понятно!"
$ws.Range("H8").Value = "Воодушевление, уверенность, спокойствие"
$ws.Range("I8").Value = "Hhh"
$ws.Range("J8").Value = "Да"
$ws.Range("K8").Value = "Нет"
$ws.Range("L8").Value = "Рост в профессии"

$ws.Range("A9").Value = "g"
$ws.Range("B9").Value = "g"
$ws.Range("C9").Value = "g"
$ws.Range("D9").Value = "g"
$ws.Range("E9").Value = "Рисовала на бесплатном марафоне😁"
$ws.Range("F9").Value = "просто получили удовольствие от процесса"
$ws.Range("G9").Value = "Всё понятно!"
$ws.Range("H9").Value = "Напряжение усилилось"
$ws.Range("I9").Value = "рпоп"
$ws.Range("J9").Value = "Нет"
$ws.Range("K9").Value = "Нет"
$ws.Range("L9").Value = "Переезд"

$ws.Range("A10").Value = "but"
$ws.Range("B10").Value = "wad"
$ws.Range("C10").Value = "wad"
$ws.Range("D10").Value = "wqedas"
$ws.Range("E10").Value = "Рисовала другой алгоритм на очном МК"
$ws.Range("F10").Value = "просто получили удовольствие от процесса"
$ws.Range("G10").Value = "Всё понятно!"
$ws.Range("H10").Value = "Хорошее настроение,вдохновение"
$ws.Range("I10").Value = "qwsA"
$ws.Range("J10").Value = "Нет"
$ws.Range("K10").Value = "Да"
$ws.Range("L10").Value = "Отношения с детьми"

foreach ($addr in $numericTextCells) {
    $ws.Range($addr).Style = "Normal"
}
